{"js": "// Add the \"Clase 12-07-2021\" section after the existing content.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Last paragraph currently in the document (the Radiogroup/radiobutton paragraph).\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Create every new paragraph first, walking forward from the existing\n// (unformatted) last paragraph, so none of them inherit character\n// formatting from one another. Character/paragraph formatting is applied\n// afterwards, only to the specific paragraphs that need it.\nconst blank1 = lastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nconst blank2 = blank1.insertParagraph(\"\", Word.InsertLocation.after);\nconst heading = blank2.insertParagraph(\"Clase 12-07-2021\", Word.InsertLocation.after);\nconst p1 = heading.insertParagraph(\n  \"En esta clase aprendimos como pasar de una activity a otra. Esto lo realizamos con la clase Intent la cual nos permite enviar informaci\u00f3n de que activity parte hacia que activity ir.\",\n  Word.InsertLocation.after\n);\nconst p2 = p1.insertParagraph(\n  \"Tambi\u00e9n vimos como pasar un valor desde una activity a otra.\",\n  Word.InsertLocation.after\n);\nconst p3 = p2.insertParagraph(\"\", Word.InsertLocation.after);\n\n// Now apply formatting \u2014 only to the paragraphs that need it.\nheading.font.bold = true;\nheading.font.italic = true;\nheading.font.underline = Word.UnderlineType.single;\n\np3.alignment = Word.Alignment.centered;\n\nawait context.sync();\n", "ps1": "# Add the \"Clase 12-07-2021\" section to the end of the document:\n#   - two blank paragraphs\n#   - a bold/italic/underlined heading paragraph (\"Clase 12-07-2021\")\n#   - two plain body paragraphs\n#   - a trailing, empty, centered paragraph\n$d = $word.ActiveDocument\n\n# Collapse a range sitting at the very end of the document body (just\n# after the last existing paragraph's text, before sectPr) so the new\n# paragraphs land after everything that is already there.\n$endRange = $d.Content\n$endRange.Collapse(0)\n\n$xml = @'\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p/><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val=\"single\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val=\"single\"/></w:rPr><w:t>Clase 12-07-2021</w:t></w:r></w:p><w:p><w:r><w:t>En esta clase aprendimos como pasar de una activity a otra. Esto lo realizamos con la clase Intent la cual nos permite enviar informaci\u00f3n de que activity parte hacia que activity ir.</w:t></w:r></w:p><w:p><w:r><w:t>Tambi\u00e9n vimos como pasar un valor desde una activity a otra.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n$endRange.InsertXML($xml)\n"}
